# Update "data-siswa" worksheet: laporan keuangan / sms gateway table
# Column A: NIM (mix of a text ID for row 2, numeric phone/ids for rows 3-7)
# Column B: SISWA (student / contact names)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - NIM stays a text value ("081840" looks numeric, so force text
# storage the same way Excel does for numbers typed into a text-formatted
# cell, to preserve the leading zero) / SISWA name update
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "081840"
$ws.Range("B2").Value = "AYAH UJEB"

# Row 3 - numeric id / name
$ws.Range("A3").Value = 24253
$ws.Range("B3").Value = "Rina"

# Row 4 - numeric id / name
$ws.Range("A4").Value = 24524525
$ws.Range("B4").Value = "Muhammad Khairu Mubarak Huzaifah"

# Row 5 - numeric id / name
$ws.Range("A5").Value = 3425115
$ws.Range("B5").Value = "Abang Ujeb keren"

# Row 6 - numeric id / name
$ws.Range("A6").Value = 64264647
$ws.Range("B6").Value = "MUHAMMAD HUZAIFAH, S.Kom."

# Row 7 - new row added at the end of the table
$ws.Range("A7").Value = 67676647
$ws.Range("B7").Value = "ira iru"
